# Registeradeal Page Locators.xlsx
# Adds a new "Locator Type" column (D) to Sheet1, marking every existing
# locator row as "CSS" except the last row (Xpath-based cookie-banner
# locator), matching the "Execute Tests on Android mobile browser" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: D1 = "Locator Type" ---------------------------------
$ws.Range("D1").Value = "Locator Type"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Size = 13
$ws.Range("D1").Font.Name = "Helvetica Neue"

# --- Data rows: D2:D44 = "CSS" (all locators use CSS selectors) -----------
$ws.Range("D2:D44").Value = "CSS"

# --- Last row uses an XPath locator instead --------------------------------
$ws.Range("D45").Value = "Xpath"

# Resize the new column to fit its contents, like Excel does automatically
# when a column's data is entered/edited.
$ws.Columns("D").AutoFit() | Out-Null

# Leave the selection where the author left off after typing the last value.
$ws.Range("D46").Select() | Out-Null
